# Update 2020 (column P) input values on the "Inputs" sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inputs")

$ws.Range("P2").Value = 28841
$ws.Range("P2").ClearFormats()
$ws.Range("P3").Value = 4973
$ws.Range("P3").ClearFormats()
$ws.Range("P4").Value = 6605
$ws.Range("P4").ClearFormats()
$ws.Range("P5").Value = 2002
$ws.Range("P5").ClearFormats()
$ws.Range("P6").Value = 6522
$ws.Range("P6").ClearFormats()
$ws.Range("P7").Value = 1539
$ws.Range("P7").ClearFormats()
$ws.Range("P8").Value = 21348
$ws.Range("P8").ClearFormats()
$ws.Range("P9").Value = 4535
$ws.Range("P9").ClearFormats()
$ws.Range("P10").Value = 4419
$ws.Range("P10").ClearFormats()
$ws.Range("P11").Value = 10962
$ws.Range("P11").ClearFormats()
$ws.Range("P12").Value = 13124
$ws.Range("P12").ClearFormats()
$ws.Range("P13").Value = 6617
$ws.Range("P13").ClearFormats()
$ws.Range("P14").Value = 5487
$ws.Range("P14").ClearFormats()
$ws.Range("P15").Value = 43699
$ws.Range("P15").ClearFormats()
$ws.Range("P16").Value = 25267
$ws.Range("P16").ClearFormats()
$ws.Range("P17").Value = 17293
$ws.Range("P17").ClearFormats()
$ws.Range("P18").Value = 10926
$ws.Range("P18").ClearFormats()
$ws.Range("P19").Value = 9503
$ws.Range("P19").ClearFormats()
$ws.Range("P20").Value = 7676
$ws.Range("P20").ClearFormats()
$ws.Range("P22").Value = 2142
$ws.Range("P23").Value = 464
$ws.Range("P26").Value = 11977
$ws.Range("P26").ClearFormats()
$ws.Range("P28").Value = 4
$ws.Range("P28").ClearFormats()
$ws.Range("P29").Value = 4
$ws.Range("P30").Value = 11981
$ws.Range("P30").ClearFormats()
$ws.Range("P33").Value = 192975
$ws.Range("P33").ClearFormats()
$ws.Range("P34").Value = 7988
$ws.Range("P35").Value = 198792
$ws.Range("P35").ClearFormats()
$ws.Range("P36").Value = 458046.08294930874
$ws.Range("P36").NumberFormat = "0"
$ws.Range("P37").Value = 0.03354211505047582
$ws.Range("P38").Value = 13905.776600000012
$ws.Range("P38").ClearFormats()
$ws.Range("P39").Value = 901
$ws.Range("P39").ClearFormats()
$ws.Range("P40").Value = 75053
$ws.Range("P41").Value = 28533
$ws.Range("P41").Interior.Color = 15773696
$ws.Range("P41").Interior.Pattern = 1
$ws.Range("P42").Value = 3737
$ws.Range("P43").Value = 1199
